# Mancala game_params.xlsx edit:
# Re-sort/move the "Sow" section rows (was already sorted by A/I/H) and
# filter the sheet down to the "Sow" tab rows, following on from moving
# visit_opp/mlap_cont to the 2nd column (row 27). The H column counters
# for the 3 "End Game" rows (76-78) shift down by one since a row was
# removed from that counting sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 (quitter / Quit Game Seeds) lost its distinct highlight style
#     now that visit_opp/mlap_cont moved to the 2nd column ahead of it ---
$ws.Range("A27:I27").Style = "Normal"

# --- End Game seed row/col counters shift down by one (9/10/11 -> 8/9/10) ---
$ws.Range("H76").Value = 8
$ws.Range("H77").Value = 9
$ws.Range("H78").Value = 10

# --- Filter column A ("tab") down to just the "Sow" rows ---
$rng = $ws.Range("A1:I88")
$rng.AutoFilter(1, @("Sow"))

# --- Last user selection before saving ---
$ws.Range("H79").Select()
